$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "Git"
$ws.Range("B8").Value = "Git Development Team"
$ws.Range("C8").Value = "2.44.0"
$ws.Range("D8").Value = "N/A"
$ws.Range("E8").Value = "GPL-2.0-or-later"
$ws.Range("F8").Value = "Perpetual"
$ws.Range("G8").Value = "Development"
$ws.Range("H8").Value = "Used for source control management"

$ws.Range("H8").Select() | Out-Null
